# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the associated timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 12:41:19"

# --- zh-cn sheet ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 12:39:11"

# --- de-de sheet ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 12:41:19"

# Widen the Status columns to fit the new, longer text (closest value
# Excel's quantized ColumnWidth property can represent).
$wsOverview.Range("E1:F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
